$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-coerced to a number
# (losing trailing zeros / exact formatting) are pre-formatted as Text so
# the literal string is preserved exactly, matching the source data feed.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D14", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.621.57"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.010.73"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "510.90"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "139.60"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "7.56"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "0.366"
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("D12").Value = "3.525.34"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "0.0000164"
$ws.Range("E15").Value = "  +7.44%  "
$ws.Range("D16").Value = "57.599.74"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "6.21"
$ws.Range("E17").Value = "  +6.19%  "
$ws.Range("D18").Value = "3.013.12"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  +3.97%  "
$ws.Range("D20").Value = "7.98"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "330.90"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "0.499"
$ws.Range("E23").Value = "  +4.54%  "
$ws.Range("D24").Value = "64.59"
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "0.0₃0924"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("D28").Value = "6.81"
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +5.17%  "
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("D32").Value = "20.61"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "4.74"
$ws.Range("E33").Value = "  +5.87%  "
$ws.Range("D34").Value = "154.76"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").Value = "5.88"
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").Value = "24.46"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").Value = "0.0677"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "3.045.61"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "37.46"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  +7.33%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "2.237.99"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "1.41"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").Value = "0.986"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "6.03"
$ws.Range("E47").Value = "  +4.99%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "19.43"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "1.87"
$ws.Range("E50").Value = "  -5.53%  "
$ws.Range("D51").Value = "0.0894"
$ws.Range("E51").Value = "  +2.87%  "
